$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.803.25"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "2.233.90"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.81"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.07%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.407"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("E11").Value = "  +4.51%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "2.565.93"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.807"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "2.249.06"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("D19").Value = "41.775.73"
$ws.Range("E19").Value = "  +2.60%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0908"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.143"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("E32").Value = "  -5.99%  "
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.18%  "
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0656"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("E37").Value = "  -8.22%  "
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000239"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0956"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.13%  "
$ws.Range("D48").Value = "1.472.27"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.71%  "
$ws.Range("E51").Value = "  -3.24%  "
